# checkin_test.docx edit:
# - merge "Name1".."Name5" paragraphs into one, prefixed by "Vorname und
#   Nachname:" and separated by line breaks
# - insert a brand-new "Namen mit Zusatz falls jünger als 6" paragraph with
#   placeholders nameTestpflicht1..5 (and the relocated _GoBack bookmark)
# - merge "Anzahl Schlüssel" / "Schlüsselkennungen" paragraphs into one,
#   separated by a line break
# - tweak the "Testdaten" paragraph wording and collapse the testdatum2..7
#   label/colon runs into single runs

$d = $word.ActiveDocument
$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

function Merge-NextParagraphs($doc, $index, $count) {
    # Deletes the paragraph mark ending paragraph $index, $count times,
    # folding that many following paragraphs into it (no break inserted).
    $p = $doc.Paragraphs.Item($index)
    for ($n = 0; $n -lt $count; $n++) {
        $pilcrow = $doc.Range($p.Range.End - 1, $p.Range.End)
        $pilcrow.Delete()
    }
}

# ---------------------------------------------------------------------
# 1) "Name1: {name1}" .. "Name5: {name5}" -> one paragraph
# ---------------------------------------------------------------------
$idx = Get-ParaIndexByText $d "Name1"
Merge-NextParagraphs $d $idx 4
$p = $d.Paragraphs.Item($idx)
$xml = "<w:p $wdNS>" +
  "<w:r><w:t>Vorname und Nachname:</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>Name1: {name1}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>Name2: {name2}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>Name3: {name3}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>Name4: {name4}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>Name5: {name5}</w:t></w:r>" +
  "</w:p>"
[void]$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 2) brand-new "Namen mit Zusatz falls jünger als 6" paragraph, inserted
#    right after the names paragraph
# ---------------------------------------------------------------------
$idx = Get-ParaIndexByText $d "Vorname und Nachname"
$p = $d.Paragraphs.Item($idx)
[void]$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($idx + 1)
$xml = "<w:p $wdNS>" +
  "<w:r><w:t>Namen mit Zusatz falls j" + [char]0xFC + "nger als 6</w:t></w:r>" +
  "<w:r><w:br/><w:t>N</w:t></w:r>" +
  "<w:r><w:t>ameTestpflicht1</w:t></w:r>" +
  "<w:r><w:t>: {nameTestpflicht1}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>NameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>2</w:t></w:r>" +
  "<w:r><w:t>: {nameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>2</w:t></w:r>" +
  "<w:r><w:t>}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>NameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>3</w:t></w:r>" +
  "<w:r><w:t>: {nameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>3</w:t></w:r>" +
  "<w:r><w:t>}</w:t></w:r>" +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>NameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>4</w:t></w:r>" +
  "<w:r><w:t>: {nameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>4</w:t></w:r>" +
  "<w:r><w:t>}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>NameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>5</w:t></w:r>" +
  "<w:r><w:t>: {nameTestpflicht</w:t></w:r>" +
  "<w:r><w:t>5</w:t></w:r>" +
  "<w:r><w:t>}</w:t></w:r>" +
  "</w:p>"
[void]$p2.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) "Anzahl Schlüssel" / "Schlüsselkennungen" -> one paragraph
# ---------------------------------------------------------------------
$idx = Get-ParaIndexByText $d "Anzahl Schl"
Merge-NextParagraphs $d $idx 1
$p = $d.Paragraphs.Item($idx)
$uum = [char]0xFC
$xml = "<w:p $wdNS>" +
  "<w:r><w:t>Anzahl Schl" + $uum + "ssel: {</w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r><w:t>anzahlSchluessel</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r><w:t>}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:t>S</w:t></w:r>" +
  "<w:r><w:t>chl" + $uum + "</w:t></w:r>" +
  "<w:r><w:t>ssel</w:t></w:r>" +
  "<w:r><w:t>kennungen</w:t></w:r>" +
  "<w:r><w:t>: {</w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r><w:t>schluessel</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r><w:t>}</w:t></w:r>" +
  "</w:p>"
[void]$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) "Testdaten" paragraph: reword intro, collapse testdatumN label+colon
# ---------------------------------------------------------------------
$idx = Get-ParaIndexByText $d "Testdaten"
$p = $d.Paragraphs.Item($idx)
$xml = "<w:p $wdNS>" +
  "<w:r><w:t>Testdaten</w:t></w:r>" +
  '<w:r><w:t xml:space="preserve"> (jeweils im Abstand von 3 Tagen)</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">testdatum2: </w:t></w:r>' +
  "<w:r><w:t>{testdatum2}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  '<w:r><w:t xml:space="preserve">testdatum3: </w:t></w:r>' +
  "<w:r><w:t>{testdatum3}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  '<w:r><w:t xml:space="preserve">testdatum4: </w:t></w:r>' +
  "<w:r><w:t>{testdatum4}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  '<w:r><w:t xml:space="preserve">testdatum5: </w:t></w:r>' +
  "<w:r><w:t>{testdatum5}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  '<w:r><w:t xml:space="preserve">testdatum6: </w:t></w:r>' +
  "<w:r><w:t>{testdatum6}</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  '<w:r><w:t xml:space="preserve">testdatum7: </w:t></w:r>' +
  "<w:r><w:t>{testdatum7}</w:t></w:r>" +
  "</w:p>"
[void]$p.Range.InsertXML($xml)

Write-Host "done"
